$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title run: split "SEVES - Alimentaire - Evenement produit" so the
#    second half instead reads "Investigation cas humain" and is bold,
#    while the first half ("SEVES - Alimentaire - ") keeps its original
#    (non-bold) formatting.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

$replaceRange = $d.Range($titleRange.Start, $titleRange.End)
$found = $replaceRange.Find.Execute("Évènement produit", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $replaceRange.Text = "Investigation cas humain"
    $replaceRange.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# 2) Move the "{%p endfor %}" text currently living in the paragraph that
#    carries the (portrait) sectPr into a brand-new paragraph right before
#    it, so the section-break paragraph itself is left empty. Previously the
#    text rode along on the section-break paragraph, which pushed it onto
#    the following (landscape) page; moving it out fixes the page
#    orientation for that line.
# ---------------------------------------------------------------------------
$cr = [char]13
$ff = [char]12

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    $trimmed = $text.TrimEnd($cr, $ff)
    $isSectionBreak = $text.EndsWith([string]$ff)
    if ($trimmed -eq "{%p endfor %}" -and $isSectionBreak) {
        $targetIndex = $i
    }
}

if ($targetIndex -ge 2) {
    $prevPara = $d.Paragraphs.Item($targetIndex - 1)
    $prevPara.Range.InsertParagraphAfter() | Out-Null

    # The freshly-inserted (empty) paragraph now occupies the slot the
    # section-break paragraph used to be at; the section-break paragraph
    # (still holding the old text + the sectPr) was pushed one slot later.
    $newTextPara = $d.Paragraphs.Item($targetIndex)
    $sectionBreakPara = $d.Paragraphs.Item($targetIndex + 1)

    $newTextPara.Range.Text = "{%p endfor %}"
    $sectionBreakPara.Range.Text = ""
}
